$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1721.5714
$ws.Range("J29").Value = 2965
$ws.Range("L29").Value = 8895
$ws.Range("N29").Value = -9457

# Row 38
$ws.Range("H38").Value = 641.36365
$ws.Range("I38").Value = 158.28572
$ws.Range("J38").Value = 1486.75
$ws.Range("K38").Value = 474.85716
$ws.Range("L38").Value = 4460.25
$ws.Range("M38").Value = -102.85716
$ws.Range("N38").Value = -5204.25

# Row 58
$ws.Range("H58").Value = 2523
$ws.Range("J58").Value = 2623.5
$ws.Range("L58").Value = 7870.5
$ws.Range("N58").Value = -8170.5

# Row 116
$ws.Range("H116").Value = 4038.6667
$ws.Range("I116").Value = 3944.4443
$ws.Range("K116").Value = 3944.4443
$ws.Range("M116").Value = -502.4443000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1269.1
$ws.Range("I45").Value = 1211.375
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1211.375
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -834.375
$ws.Range("N45").Value = -2254

# Row 110
$ws.Range("H110").Value = 7825.35
$ws.Range("I110").Value = 5965.4116
$ws.Range("K110").Value = 5965.4116
$ws.Range("M110").Value = -3920.4116

# Row 122
$ws.Range("H122").Value = 3069.9333
$ws.Range("I122").Value = 1879.0834
$ws.Range("K122").Value = 5637.2502
$ws.Range("M122").Value = -3187.2502

# Row 133
$ws.Range("H133").Value = 74122.625
$ws.Range("J133").Value = 74122.625
$ws.Range("L133").Value = 74122.625
$ws.Range("N133").Value = -79182.625

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 28315.578
$ws.Range("I20").Value = 19492.875
$ws.Range("J20").Value = 38398.668
$ws.Range("K20").Value = 19492.875
$ws.Range("L20").Value = 38398.668
$ws.Range("M20").Value = -19245.875
$ws.Range("N20").Value = -38892.668

# Row 26
$ws.Range("H26").Value = 75000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 75000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 75000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -75584

# Row 96
$ws.Range("H96").Value = 37476
$ws.Range("I96").Value = 12428
$ws.Range("K96").Value = 12428
$ws.Range("M96").Value = -9682

# Row 107
$ws.Range("H107").Value = 2206.25
$ws.Range("I107").Value = 1863.68
$ws.Range("K107").Value = 1863.68
$ws.Range("M107").Value = 56.31999999999994

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1436133.9
$ws.Range("J5").Value = 4067223.2
$ws.Range("L5").Value = 12201669.6
$ws.Range("N5").Value = -12201893.6

# Row 86
$ws.Range("H86").Value = 658.3939
$ws.Range("I86").Value = 671.6957
$ws.Range("J86").Value = 627.8
$ws.Range("K86").Value = 2015.0871
$ws.Range("L86").Value = 1883.4
$ws.Range("M86").Value = -829.0871
$ws.Range("N86").Value = -4255.4

# Row 89
$ws.Range("H89").Value = 658.3939
$ws.Range("I89").Value = 671.6957
$ws.Range("J89").Value = 627.8
$ws.Range("K89").Value = 6045.2613
$ws.Range("L89").Value = 5650.2
$ws.Range("M89").Value = -117.2613000000001
$ws.Range("N89").Value = -17506.2

# Row 100
$ws.Range("H100").Value = 4312
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4312
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 12936
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -14558

# Row 114
$ws.Range("H114").Value = 2825.1428
$ws.Range("I114").Value = 1892.6666
$ws.Range("J114").Value = 3524.5
$ws.Range("K114").Value = 5677.9998
$ws.Range("L114").Value = 10573.5
$ws.Range("M114").Value = -2423.9998
$ws.Range("N114").Value = -17081.5

# Row 128
$ws.Range("H128").Value = 188220.75
$ws.Range("I128").Value = 188220.75
$ws.Range("K128").Value = 564662.25
$ws.Range("M128").Value = -559682.25

# Row 135
$ws.Range("H135").Value = 1436133.9
$ws.Range("J135").Value = 4067223.2
$ws.Range("L135").Value = 36605008.8
$ws.Range("N135").Value = -36610078.8

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4141.857
$ws.Range("I102").Value = 4998.5
$ws.Range("K102").Value = 4998.5
$ws.Range("M102").Value = -3376.5

# Row 103
$ws.Range("H103").Value = 300000
$ws.Range("J103").Value = 300000
$ws.Range("L103").Value = 300000
$ws.Range("N103").Value = -302344

# Row 122
$ws.Range("H122").Value = 2484.7368
$ws.Range("I122").Value = 1737
$ws.Range("K122").Value = 5211
$ws.Range("M122").Value = -2761

# Row 126
$ws.Range("H126").Value = 4740.4
$ws.Range("I126").Value = 3694.0667
$ws.Range("K126").Value = 11082.2001
$ws.Range("M126").Value = -8612.2001

# Row 134
$ws.Range("H134").Value = 17860.4
$ws.Range("J134").Value = 17860.4
$ws.Range("L134").Value = 53581.2
$ws.Range("N134").Value = -58651.2

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 35999
$ws.Range("J38").Value = 35999
$ws.Range("L38").Value = 35999
$ws.Range("N38").Value = -36819

# Row 61
$ws.Range("H61").Value = 4236.3184
$ws.Range("I61").Value = 2393.9
$ws.Range("K61").Value = 2393.9
$ws.Range("M61").Value = -2191.9

# Row 113
$ws.Range("H113").Value = 4236.3184
$ws.Range("I113").Value = 2393.9
$ws.Range("K113").Value = 2393.9
$ws.Range("M113").Value = -223.9000000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# Row 81
$ws.Range("H81").Value = 1303.2142
$ws.Range("I81").Value = 945.4167
$ws.Range("J81").Value = 3450
$ws.Range("K81").Value = 1890.8334
$ws.Range("L81").Value = 6900
$ws.Range("M81").Value = -829.8334
$ws.Range("N81").Value = -9022

# Row 84
$ws.Range("H84").Value = 1303.2142
$ws.Range("I84").Value = 945.4167
$ws.Range("J84").Value = 3450
$ws.Range("K84").Value = 9454.166999999999
$ws.Range("L84").Value = 34500
$ws.Range("M84").Value = -4150.166999999999
$ws.Range("N84").Value = -45108

# Row 96
$ws.Range("H96").Value = 1989.3889
$ws.Range("I96").Value = 1615.5714
$ws.Range("J96").Value = 2227.2727
$ws.Range("K96").Value = 1615.5714
$ws.Range("L96").Value = 2227.2727
$ws.Range("M96").Value = -242.5714
$ws.Range("N96").Value = -4973.2727

# Row 136
$ws.Range("H136").Value = 7773.073
$ws.Range("I136").Value = 1956.2333
$ws.Range("K136").Value = 5868.699900000001
$ws.Range("M136").Value = -3318.699900000001
